# Update the "Inventaire" sheet: record the two notice cells and leave the
# selection where the author left it (F7) when they saved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "mon code est : "
$ws.Range("C3").Value = "t'as vraimnet cru que j'allais te le donner"

[void]$ws.Range("F7").Select()
